# Daily attendance processing - reorders the "Recorded By" (column G) names
# so that "System" / "system" entries are listed first, followed by the
# human recorder(s), for every session row where the recorded-by list
# currently ends with "System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Reverse-CommaList {
    param([string]$text)
    $parts = $text -split ", "
    $count = $parts.Count
    $reversed = @()
    for ($i = $count - 1; $i -ge 0; $i--) {
        $reversed += $parts[$i]
    }
    return [string]::Join(", ", $reversed)
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

# "Recorded By" lives in column G (7); fall back to scanning the header
# row in case the sheet layout ever shifts.
$recordedByCol = 7
for ($c = 1; $c -le $lastCol; $c++) {
    if ($ws.Cells.Item(1, $c).Text -eq "Recorded By") {
        $recordedByCol = $c
        break
    }
}

$changedCount = 0
for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $recordedByCol)
    $value = $cell.Text

    if ([string]::IsNullOrEmpty($value)) {
        continue
    }

    if ($value.ToLower().EndsWith("system") -and -not $value.Contains("admin@admin.com")) {
        $newValue = Reverse-CommaList $value
        if ($newValue -ne $value) {
            $cell.Value = $newValue
            $changedCount++
        }
    }
}

Write-Host "Reordered Recorded By entries in $changedCount row(s)."
